$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E24) currently lists periods in descending
# order (2103 down to 2007). The new batch of account-statement data being
# appended requires the list to run in ascending order (2007 up to 2103),
# keeping the special "Valor Mora" amount (30666) attached to period 2103
# (now the last row) while every other period keeps the standard 40000.

$periods = @("2007", "2008", "2009", "2010", "2011", "2012", "2101", "2102", "2103")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $period = $periods[$i]

    $ws.Range("E$row").Value = $period

    if ($period -eq "2103") {
        $ws.Range("F$row").Value = 30666
    } else {
        $ws.Range("F$row").Value = 40000
    }
}
